$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 21.487739
$ws.Range("C3").Value = 33.449434
$ws.Range("D3").Value = 36.351233
$ws.Range("E3").Value = 39.081447
$ws.Range("F3").Value = 42.737567
# Row 4
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 0.24
$ws.Range("C4").Value = 0.25
$ws.Range("D4").Value = 0.08
$ws.Range("E4").Value = 0.13
$ws.Range("F4").Value = 0.04
# Row 5
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 0.05
$ws.Range("C5").Value = 0.03
$ws.Range("D5").Value = 0.2
$ws.Range("E5").Value = 0.07000000000000001
$ws.Range("F5").Value = 0.02
# Row 6
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 0.037014
$ws.Range("C6").Value = 0.030728
$ws.Range("D6").Value = 0.010464
$ws.Range("E6").Value = 0.010177
$ws.Range("F6").Value = 0.01817
# Row 7
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 0.132336
$ws.Range("C7").Value = 0.06342
$ws.Range("D7").Value = 0.248894
$ws.Range("E7").Value = 0.46513
$ws.Range("F7").Value = 0.494593
# Row 8
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = 0.083643
$ws.Range("C8").Value = 0.08955200000000001
$ws.Range("D8").Value = 0.155616
$ws.Range("E8").Value = 0.054392
$ws.Range("F8").Value = 0.10753
# Row 9
$ws.Range("A9").Value = 0
$ws.Range("B9").Value = 0.468956
$ws.Range("C9").Value = 0.518135
$ws.Range("D9").Value = 0.198485
$ws.Range("E9").Value = 0.015336
$ws.Range("F9").Value = 0.013869
# Row 10
$ws.Range("A10").Value = 0
$ws.Range("B10").Value = 0.086788
$ws.Range("C10").Value = 0.082331
$ws.Range("D10").Value = 0.051939
$ws.Range("E10").Value = 0.011576
$ws.Range("F10").Value = 0.014201
# Row 11
$ws.Range("A11").Value = 0
$ws.Range("B11").Value = 0.132094
$ws.Range("C11").Value = 0.059683
$ws.Range("D11").Value = 0.282421
$ws.Range("E11").Value = 0.335834
$ws.Range("F11").Value = 0.263771
# Row 12
$ws.Range("A12").Value = 0
$ws.Range("B12").Value = 0.05917
$ws.Range("C12").Value = 0.156151
$ws.Range("D12").Value = 0.052179
$ws.Range("E12").Value = 0.107554
$ws.Range("F12").Value = 0.087865
# Row 13
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 0.758576
$ws.Range("C13").Value = 0.834748
$ws.Range("D13").Value = 0.417218
$ws.Range("E13").Value = 0.152476
$ws.Range("F13").Value = 0.074265
# Row 14
$ws.Range("A14").Value = 0
$ws.Range("B14").Value = 0.082551
$ws.Range("C14").Value = 0.053798
$ws.Range("D14").Value = 0.096844
$ws.Range("E14").Value = 0.359644
$ws.Range("F14").Value = 0.192006
# Row 15
$ws.Range("A15").Value = 0
$ws.Range("B15").Value = 0.017988
$ws.Range("C15").Value = 0.006497
$ws.Range("D15").Value = 0.281691
$ws.Range("E15").Value = 0.287962
$ws.Range("F15").Value = 0.556658
# Row 16
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = 0.01245
$ws.Range("C16").Value = 0.016658
$ws.Range("D16").Value = 0.023613
$ws.Range("E16").Value = 0.015358
$ws.Range("F16").Value = 0.007247
# Row 17
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = 0.123306
$ws.Range("C17").Value = 0.081051
$ws.Range("D17").Value = 0.029601
$ws.Range("E17").Value = 0.129525
$ws.Range("F17").Value = 0.114049
# Row 18
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = 0.005129
$ws.Range("C18").Value = 0.007248
$ws.Range("D18").Value = 0.151031
$ws.Range("E18").Value = 0.055035
$ws.Range("F18").Value = 0.055776
# Row 19
$ws.Range("A19").Value = 0
$ws.Range("B19").Value = 0.007336
$ws.Range("C19").Value = 0.038562
$ws.Range("D19").Value = 0.00651
$ws.Range("E19").Value = 0.00468
$ws.Range("F19").Value = 0.00311
# Row 20
$ws.Range("A20").Value = 0
$ws.Range("B20").Value = 0.873723
$ws.Range("C20").Value = 0.765254
$ws.Range("D20").Value = 0.757493
$ws.Range("E20").Value = 0.547498
$ws.Range("F20").Value = 0.897111
# Row 21
$ws.Range("A21").Value = 0
$ws.Range("B21").Value = 0.118941
$ws.Range("C21").Value = 0.196184
$ws.Range("D21").Value = 0.235997
$ws.Range("E21").Value = 0.447822
$ws.Range("F21").Value = 0.09977999999999999
# Row 22
$ws.Range("A22").Value = 0.06743775
# Row 23
$ws.Range("A23").Value = 0.29886375
# Row 24
$ws.Range("A24").Value = 0.06159575
# Row 25
$ws.Range("A25").Value = 0.30633225
# Row 26
$ws.Range("A26").Value = 0.09225925
# Row 27
$ws.Range("A27").Value = 0.1257095
# Row 28
$ws.Range("A28").Value = 0.04780175
# Row 29
$ws.Range("A29").Value = 0.06743775
$ws.Range("B29").Value = 0.29886375
$ws.Range("C29").Value = 0.06159575
$ws.Range("D29").Value = 0.30633225
$ws.Range("E29").Value = 0.09225925
$ws.Range("F29").Value = 0.1257095
$ws.Range("G29").Value = 0.04780175
# Row 30
$ws.Range("A30").Value = 0.5911573631587611
# Row 31
$ws.Range("A31").Value = 0.006152395
# Row 32
$ws.Range("A32").Value = 0.306738115
# Row 33
$ws.Range("A33").Value = 0.23492557
# Row 34
$ws.Range("A34").Value = 0.179188181
# Row 35
$ws.Range("A35").Value = 0.272995739
# Row 36
$ws.Range("A36").Value = 0.5911573631587611
$ws.Range("B36").Value = 0.006152395
$ws.Range("C36").Value = 0.306738115
$ws.Range("D36").Value = 0.23492557
$ws.Range("E36").Value = 0.179188181
$ws.Range("F36").Value = 0.272995739
# Row 37
$ws.Range("A37").Value = 0.139625
# Row 38
$ws.Range("A38").Value = 0.08347500000000002
# Row 39
$ws.Range("A39").Value = 0.018777915
# Row 40
$ws.Range("A40").Value = 0.29610506375
# Row 41
$ws.Range("A41").Value = 0.098933745
# Row 42
$ws.Range("A42").Value = 0.215002365
# Row 43
$ws.Range("A43").Value = 0.0449965425
# Row 44
$ws.Range("A44").Value = 0.2306732825
# Row 45
$ws.Range("A45").Value = 0.0955102625
# Row 46
$ws.Range("A46").Value = 0.018777915
$ws.Range("B46").Value = 0.29610506375
$ws.Range("C46").Value = 0.098933745
$ws.Range("D46").Value = 0.215002365
$ws.Range("E46").Value = 0.0449965425
$ws.Range("F46").Value = 0.2306732825
$ws.Range("G46").Value = 0.0955102625
# Row 47
$ws.Range("A47").Value = 0.41376086375
# Row 48
$ws.Range("A48").Value = 0.1723875250000001
# Row 49
$ws.Range("A49").Value = 0.24294294375
# Row 50
$ws.Range("A50").Value = 0.0160233675
# Row 51
$ws.Range("A51").Value = 0.0912479375
# Row 52
$ws.Range("A52").Value = 0.06363703125000002
# Row 53
$ws.Range("A53").Value = 0.41376086375
$ws.Range("B53").Value = 0.1723875250000001
$ws.Range("C53").Value = 0.24294294375
$ws.Range("D53").Value = 0.0160233675
$ws.Range("E53").Value = 0.0912479375
$ws.Range("F53").Value = 0.06363703125000002
# Row 54
$ws.Range("A54").Value = 0.01227494
# Row 55
$ws.Range("A55").Value = 0.7386553650000001
# Row 56
$ws.Range("A56").Value = 0.24906986125
# Row 57
$ws.Range("A57").Value = 0.6552198814457406
# Row 58
$ws.Range("A58").Value = 0.7137426433807432
# Row 59
$ws.Range("A59").Value = 0.6844812624132419
# Row 60
$ws.Range("A60").Value = 0.618397460625
